$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 14 ("Forste iterasjon"): content placeholder "GUI" / "Klasser"
# gets two extra paragraphs appended: "Db" and a trailing empty one.
# -----------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(3)
$tr14 = $sh14.TextFrame.TextRange
$lastPara14 = $tr14.Paragraphs($tr14.Paragraphs().Count, 1)
$lastPara14.InsertAfter("`rDb`r")

# -----------------------------------------------------------------
# Slide 4 ("Scrum"): "Daily meetings" paragraph becomes
# "Stand-up meetings" (two runs: "Stand-up " / "meetings").
# -----------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange

# Locate the "Daily meetings" paragraph (2nd paragraph in the box).
$para2 = $tr4.Paragraphs(2, 1)

# Replace "Daily" (first 5 chars of the paragraph) with nothing, then
# type "Stand-up" fresh right before what remains (" meetings"), so a
# brand-new run is created instead of reusing the old "Daily" run.
$para2.Characters(1, 5).Delete()
$remainder = $tr4.Paragraphs(2, 1)
$remainder.InsertBefore("Stand-up")

# Merge the freshly typed "Stand-up" together with the single space
# that used to separate "Daily" and "meetings" into one run reading
# "Stand-up ".
$para2again = $tr4.Paragraphs(2, 1)
$merged = $para2again.Characters(1, 9)
$merged.Text = "Stand-up "

# Finally, rebuild the "meetings" run fresh as well so it no longer
# carries the old spell-check "err" marker.
$para2final = $tr4.Paragraphs(2, 1)
$meetingsOld = $para2final.Characters(10, 8)
$meetingsOld.Delete()
$tail = $tr4.Paragraphs(2, 1)
$tail.InsertAfter("meetings")
